$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 (Marking)
$ws.Range("B11").Value = 9
$ws.Range("C11").Value = 2

# Row 12 (Total)
$ws.Range("B12").Value = 207
$ws.Range("C12").Value = -6
$ws.Range("E12").Value = "201/252"
